$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: add quantity and first resistor value next to existing "resistors" label
$ws.Range("D9").Value = 7
$ws.Range("F9").Value = "1k"

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("F10").Value = "1M"

# Row 11
$ws.Range("D11").Value = 1
$ws.Range("F11").Value = "250k"

# Row 12
$ws.Range("D12").Value = 1
$ws.Range("F12").Value = "10k"

# Row 13
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = "2.5k"

# Row 14: pots
$ws.Range("A14").Value = "pots"
$ws.Range("D14").Value = 2

# Row 15 (set F15's string value before F14's so the shared-string table
# ordering matches the source ordering: "1k to 0" precedes
# "some wide range around 100k, idk")
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = "1k to 0"

$ws.Range("F14").Value = "some wide range around 100k, idk"

# Row 16: switch
$ws.Range("A16").Value = "switch"
$ws.Range("F16").Value = "preferably connects all three voltages or disconnects the 120V wall power"

# Row 17: NPN
$ws.Range("A17").Value = "NPN"

# Row 18: PNP
$ws.Range("A18").Value = "PNP"

# Update active selection to match post-edit state
$ws.Range("A19").Select()
